$p = $ppt.ActivePresentation

# --- Slide 1: Content-only slide (Title placeholder removed) ---
$s1 = $p.Slides.Add(1, 2)
$s1.Shapes.Item(1).Delete()   # remove the Title placeholder shape
$shp1 = $s1.Shapes.Item(1)    # Content Placeholder 2
$shp1_tr = $shp1.TextFrame.TextRange
$shp1_tr.Text = 'Principale problema: stanno mettendo mano alla dataplatform aziendale e stanno rifacendo: bigger problem: ingestion e estrarre dati da sistemi proprietari. SIa saleforce che SAP se devono togliere i dati hanno grossi problemi di billing a cause di lock in. 2 opzioni per analisi: - usa i loro tool - estrai i dati tramite accrocchi.'
$shp1_tr.ParagraphFormat.Bullet.Type = 0
$shp1_tr.IndentLevel = 1
$rulerLvl = $shp1.TextFrame.Ruler.Levels.Item(1)
$rulerLvl.FirstMargin = 0
$rulerLvl.LeftMargin = 0
$shp1_p1 = $shp1_tr.InsertAfter("`r" + 'Costruire una platform dato un sistema esistente il task tosto è come estrapolare i dati dai sistemi di base e migrarli nei nuovi servizi. Fanno una specie di sampling in cui perdono tanti dati (e.g. campionamento 30min), che se devi fare un analisi comportamentale è bad. C’è molta spinta per riportare sul frontend informazioni computate in real time da una data platform, in questi casi 30min sono decisamente tanti.')
$shp1_p1.ParagraphFormat.Bullet.Type = 0
$shp1_p1.IndentLevel = 1
$rulerLvl = $shp1.TextFrame.Ruler.Levels.Item(1)
$rulerLvl.FirstMargin = 0
$rulerLvl.LeftMargin = 0
$shp1_p2 = $shp1_tr.InsertAfter("`r" + 'Pipeline back&forth. - dati dabellari strutturati')
$shp1_p2.ParagraphFormat.Bullet.Type = 0
$shp1_p2.IndentLevel = 1
$rulerLvl = $shp1.TextFrame.Ruler.Levels.Item(1)
$rulerLvl.FirstMargin = 0
$rulerLvl.LeftMargin = 0
$shp1_p3 = $shp1_tr.InsertAfter("`r" + 'Sorgente è Salesforce, collection in pull e poi processi. Change Data Capture rispetto aquello che avevo prima. Quello che pensano di fare è hashare i campi per i quali vogliono capire se sono cambiati e confrontano gl ihash. Hanno la parte storica, e una parte live utilizzata per reporting, eventuali modelli da utilizzare su dati virtuali.')
$shp1_p3.ParagraphFormat.Bullet.Type = 0
$shp1_p3.IndentLevel = 1
$rulerLvl = $shp1.TextFrame.Ruler.Levels.Item(1)
$rulerLvl.FirstMargin = 0
$rulerLvl.LeftMargin = 0
$shp1_p4 = $shp1_tr.InsertAfter("`r" + 'Dopo CDC, tutta la serie storica sul cliente viene portato su un tool su cui viene fatto girare un modello fare AI (xgBoost), e.g. calcola CLV (custom lifetime value e altri due indici). Il risulato del modello viene rigirato sulla platform che a sua volta lo ripusha su Salesforce.')
$shp1_p4.ParagraphFormat.Bullet.Type = 0
$shp1_p4.IndentLevel = 1
$rulerLvl = $shp1.TextFrame.Ruler.Levels.Item(1)
$rulerLvl.FirstMargin = 0
$rulerLvl.LeftMargin = 0
$shp1_p5 = $shp1_tr.InsertAfter("`r" + 'Su Google Cloud Platform. ' + 'informati su Dataform' + ', SQLX. permette di fare ETL con source control. Skilled Query e BIGQuery(dove fanno la maggior parte delle trasformazioni). per AI utilizzano Vertex: un notebook schedulato. Puoi schedulare notebook da runnare ad un certo orario con macchine di certa dimensione (tipo databricks). Con Google fanno una riunione ogni due settimane (per consulenze di data platform, engine selection e cosa fare). Hanno usato per parecchio e usano AutoML, se non funzia passano a modelli custom (Model Garden (?)): modelli pretrainati che ti mettono a disposizione API per lavorare tipo serverless.')
$shp1_p5.ParagraphFormat.Bullet.Type = 0
$shp1_p5.IndentLevel = 1
$rulerLvl = $shp1.TextFrame.Ruler.Levels.Item(1)
$rulerLvl.FirstMargin = 0
$rulerLvl.LeftMargin = 0
$fullTr = $shp1.TextFrame.TextRange
$paraStart = $fullTr.Text.Length - ($shp1_p5.Text.Length) + 1
$boldRun = $fullTr.Characters($paraStart + 26, 21)
$boldRun.Font.Bold = $true

# --- Slide 2: Title + Content slide ---
$s2 = $p.Slides.Add(2, 2)
$shpTitle = $s2.Shapes.Item(1)   # Title 1
$shpBody = $s2.Shapes.Item(2)    # Content Placeholder 2
$shpTitle_tr = $shpTitle.TextFrame.TextRange
$shpTitle_tr.Text = 'ICO'
$shpTitle_tr.ParagraphFormat.Bullet.Type = 0
$shpTitle_tr.IndentLevel = 1
$rulerLvl = $shpTitle.TextFrame.Ruler.Levels.Item(1)
$rulerLvl.FirstMargin = 0
$rulerLvl.LeftMargin = 0
$shpBody_tr = $shpBody.TextFrame.TextRange
$shpBody_tr.Text = 'DWH è già nel silver layer Nel silver i dati sono volatili, 3/4 giorni di vita La parte bronze è una specie di storicizzazione della parte “effimera” e “Esterna” del viola'
$shpBody_tr.ParagraphFormat.Bullet.Type = 0
$shpBody_tr.IndentLevel = 1
$rulerLvl = $shpBody.TextFrame.Ruler.Levels.Item(1)
$rulerLvl.FirstMargin = 0
$rulerLvl.LeftMargin = 0
$shpBody_p1 = $shpBody_tr.InsertAfter("`r" + 'viola arancio - import arancio grigio - clean')
$shpBody_p1.ParagraphFormat.Bullet.Type = 0
$shpBody_p1.IndentLevel = 1
$rulerLvl = $shpBody.TextFrame.Ruler.Levels.Item(1)
$rulerLvl.FirstMargin = 0
$rulerLvl.LeftMargin = 0
$shpBody_p2 = $shpBody_tr.InsertAfter("`r" + 'Mosaic AI')
$shpBody_p2.ParagraphFormat.Bullet.Type = 0
$shpBody_p2.IndentLevel = 1
$rulerLvl = $shpBody.TextFrame.Ruler.Levels.Item(1)
$rulerLvl.FirstMargin = 0
$rulerLvl.LeftMargin = 0
$shpBody_p3 = $shpBody_tr.InsertAfter("`r" + 'Le nostre soluzioni funzionano a livell di orchestrazione, ma soprattutto su databricks e azure la chiave è il compute che decidi di utilizzare')
$shpBody_p3.ParagraphFormat.Bullet.Type = 0
$shpBody_p3.IndentLevel = 1
$rulerLvl = $shpBody.TextFrame.Ruler.Levels.Item(1)
$rulerLvl.FirstMargin = 0
$rulerLvl.LeftMargin = 0
$shpBody_p4 = $shpBody_tr.InsertAfter("`r" + 'Suddividere il processo in : ingestion - engineering - analisi - elaborazione avanzata')
$shpBody_p4.ParagraphFormat.Bullet.Type = 0
$shpBody_p4.IndentLevel = 1
$rulerLvl = $shpBody.TextFrame.Ruler.Levels.Item(1)
$rulerLvl.FirstMargin = 0
$rulerLvl.LeftMargin = 0
$shpBody_p5 = $shpBody_tr.InsertAfter("`r" + 'Ci sono tanti piccoli aspetti che non riusciamo a modellare: unity catalog, lineage tra trasformazioni')
$shpBody_p5.ParagraphFormat.Bullet.Type = 0
$shpBody_p5.IndentLevel = 1
$rulerLvl = $shpBody.TextFrame.Ruler.Levels.Item(1)
$rulerLvl.FirstMargin = 0
$rulerLvl.LeftMargin = 0
$shpBody_p6 = $shpBody_tr.InsertAfter("`r" + 'Databricks')
$shpBody_p6.ParagraphFormat.Bullet.Type = 0
$shpBody_p6.IndentLevel = 1
$rulerLvl = $shpBody.TextFrame.Ruler.Levels.Item(1)
$rulerLvl.FirstMargin = 0
$rulerLvl.LeftMargin = 0
